# Actualizacion 11 de Mayo - Manana
# Updates Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio/Blancos/Por_Blan
# figures across the three "Parcial" worksheets for several groups.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("E14").Value = 36
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 9
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 22
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 8.300000000000001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 100
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 7.2
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("I29").Value = 6.4
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("I32").Value = 7
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("I33").Value = 7
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("I34").Value = 7
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("I35").Value = 6.3
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("E36").Value = 23
$ws.Range("F36").Value = 13
$ws.Range("G36").Value = 63.89
$ws.Range("H36").Value = 36.11
$ws.Range("I36").Value = 6.8
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("E37").Value = 27
$ws.Range("F37").Value = 7
$ws.Range("G37").Value = 79.41
$ws.Range("H37").Value = 20.59
$ws.Range("I37").Value = 7.8
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("E38").Value = 19
$ws.Range("F38").Value = 7
$ws.Range("G38").Value = 73.08
$ws.Range("H38").Value = 26.92
$ws.Range("I38").Value = 6.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("E39").Value = 24
$ws.Range("F39").Value = 11
$ws.Range("G39").Value = 68.56999999999999
$ws.Range("H39").Value = 31.43
$ws.Range("I39").Value = 6.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 14
$ws.Range("G40").Value = 51.72
$ws.Range("H40").Value = 48.28
$ws.Range("I40").Value = 6.6
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("E14").Value = 36
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 8.9
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 22
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 8.4
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 8.9
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 100
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 7.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("E29").Value = 18
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 48.65
$ws.Range("H29").Value = 51.35
$ws.Range("I29").Value = 7.1
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = 45.95
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 29
$ws.Range("G31").Value = 21.62
$ws.Range("H31").Value = 78.38
$ws.Range("I31").Value = 6.8
$ws.Range("J31").Value = 28
$ws.Range("K31").Value = 75.68000000000001
$ws.Range("E32").Value = 15
$ws.Range("F32").Value = 16
$ws.Range("G32").Value = 48.39
$ws.Range("H32").Value = 51.61
$ws.Range("I32").Value = 8
$ws.Range("J32").Value = 12
$ws.Range("K32").Value = 38.71
$ws.Range("E33").Value = 12
$ws.Range("F33").Value = 9
$ws.Range("G33").Value = 57.14
$ws.Range("H33").Value = 42.86
$ws.Range("I33").Value = 7
$ws.Range("J33").Value = 8
$ws.Range("K33").Value = 38.1
$ws.Range("E34").Value = 26
$ws.Range("F34").Value = 14
$ws.Range("G34").Value = 65
$ws.Range("H34").Value = 35
$ws.Range("I34").Value = 7.2
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("E35").Value = 11
$ws.Range("F35").Value = 12
$ws.Range("G35").Value = 47.83
$ws.Range("H35").Value = 52.17
$ws.Range("I35").Value = 6.9
$ws.Range("J35").Value = 8
$ws.Range("K35").Value = 34.78
$ws.Range("E36").Value = 20
$ws.Range("F36").Value = 16
$ws.Range("G36").Value = 55.56
$ws.Range("H36").Value = 44.44
$ws.Range("I36").Value = 7.9
$ws.Range("J36").Value = 15
$ws.Range("K36").Value = 41.67
$ws.Range("E37").Value = 23
$ws.Range("F37").Value = 11
$ws.Range("G37").Value = 67.65000000000001
$ws.Range("H37").Value = 32.35
$ws.Range("I37").Value = 8.699999999999999
$ws.Range("J37").Value = 11
$ws.Range("K37").Value = 32.35
$ws.Range("E38").Value = 17
$ws.Range("F38").Value = 9
$ws.Range("G38").Value = 65.38
$ws.Range("H38").Value = 34.62
$ws.Range("I38").Value = 7.6
$ws.Range("J38").Value = 9
$ws.Range("K38").Value = 34.62
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = 22
$ws.Range("G39").Value = 37.14
$ws.Range("H39").Value = 62.86
$ws.Range("I39").Value = 7
$ws.Range("J39").Value = 22
$ws.Range("K39").Value = 62.86
$ws.Range("E40").Value = 14
$ws.Range("F40").Value = 15
$ws.Range("G40").Value = 48.28
$ws.Range("H40").Value = 51.72
$ws.Range("I40").Value = 7.8
$ws.Range("J40").Value = 14
$ws.Range("K40").Value = 48.28

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("E14").Value = 36
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 9.1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("E15").Value = 22
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 8.300000000000001
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 8.699999999999999
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("E17").Value = 21
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 100
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 7.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("I29").Value = 6.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("E32").Value = 21
$ws.Range("F32").Value = 10
$ws.Range("G32").Value = 67.73999999999999
$ws.Range("H32").Value = 32.26
$ws.Range("I32").Value = 7.1
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("I33").Value = 6.9
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("E34").Value = 26
$ws.Range("F34").Value = 14
$ws.Range("G34").Value = 65
$ws.Range("H34").Value = 35
$ws.Range("I34").Value = 7.1
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("E35").Value = 14
$ws.Range("F35").Value = 9
$ws.Range("G35").Value = 60.87
$ws.Range("H35").Value = 39.13
$ws.Range("I35").Value = 6.4
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("E36").Value = 25
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 69.44
$ws.Range("H36").Value = 30.56
$ws.Range("I36").Value = 7
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("E37").Value = 27
$ws.Range("F37").Value = 7
$ws.Range("G37").Value = 79.41
$ws.Range("H37").Value = 20.59
$ws.Range("I37").Value = 7.9
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("E38").Value = 19
$ws.Range("F38").Value = 7
$ws.Range("G38").Value = 73.08
$ws.Range("H38").Value = 26.92
$ws.Range("I38").Value = 7
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("E39").Value = 24
$ws.Range("F39").Value = 11
$ws.Range("G39").Value = 68.56999999999999
$ws.Range("H39").Value = 31.43
$ws.Range("I39").Value = 6.5
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 14
$ws.Range("G40").Value = 51.72
$ws.Range("H40").Value = 48.28
$ws.Range("I40").Value = 6.6
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
